# The "Förändrad" (Changed) date in column C was bumped by one day
# (2023-10-04 -> 2023-10-05, i.e. serial 45203 -> 45204) for every data
# row (rows 2 through 135) on the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C135").Value = 45204
